$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3889
$ws1.Range("F7").Value = 25
$ws1.Range("F11").Value = 1443
$ws1.Range("F12").Value = 255
$ws1.Range("F13").Value = 2571

# Sheet "全部类型" (fourth sheet) - update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3890
$ws4.Range("F7").Value = 25
$ws4.Range("F14").Value = 1443
$ws4.Range("F15").Value = 255
$ws4.Range("F16").Value = 2571

$wb.Save()
